# The "Vostro / Dell / DKS" enrollment-device rows (originally rows 6 and 7)
# are removed from the master device spec table; every row below shifts up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("6:7").Delete()

# Leave the cursor parked where the author left it after trimming the table.
[void]$ws.Range("E16").Select()

# Page was set up for printing (Letter/A4-style paper, portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
